$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A88").Value = "G1"
$ws.Range("B88").Value = "Test1"
$ws.Range("C88").Value = 45904
$ws.Range("C88").NumberFormat = $ws.Range("C87").NumberFormat
$ws.Range("D88").Value = 0.6584189185966455
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = -0.01

$ws.Range("A89").Value = "G2"
$ws.Range("B89").Value = "sedrftgyhuioygtfrd"
$ws.Range("C89").Value = 45904
$ws.Range("C89").NumberFormat = $ws.Range("C87").NumberFormat
$ws.Range("D89").Value = 0.6584189185966455
$ws.Range("E89").Value = 0
$ws.Range("F89").Value = -0.01
